# Updated Global_M2 for easier usage.
# Applies revised values to existing monthly M2 rows (122-218) and appends
# three new monthly rows (221-223) for Romania M2 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Revise existing values in columns C:F for the rows below ----------
$updates = @{
    122 = 261572700000
    123 = 259192800000
    124 = 258559200000
    125 = 255274400000
    131 = 263647500000
    134 = 286126300000
    135 = 283933800000
    136 = 283623100000
    137 = 280655000000
    142 = 295604600000
    143 = 295922500000
    144 = 296506400000
    146 = 314026000000
    147 = 309824800000
    148 = 311655400000
    149 = 314511600000
    150 = 318498700000
    151 = 320618400000
    152 = 321076200000
    153 = 324218500000
    154 = 329607000000
    155 = 332065200000
    156 = 336385900000
    157 = 339280200000
    158 = 350004800000
    159 = 348823700000
    160 = 352411100000
    161 = 351238500000
    162 = 354917800000
    163 = 356539500000
    164 = 362385300000
    165 = 360221900000
    170 = 381075300000
    171 = 382602200000
    172 = 384958100000
    173 = 383090000000
    182 = 422631600000
    192 = 469280100000
    194 = 487349900000
    195 = 490302300000
    196 = 496963100000
    197 = 499199700000
    206 = 564423000000
    210 = 569711700000
    212 = 569309400000
    215 = 581768700000
    218 = 603199600000
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("C${row}:F${row}").Value = $value
}

# --- 2) Append three new monthly rows (221-223) ----------------------------
$newRows = @(
    @{ Row = 221; DateSerial = 44986.45833333334; Value = 613926400000 },
    @{ Row = 222; DateSerial = 45017.45833333334; Value = 618680400000 },
    @{ Row = 223; DateSerial = 45047.41666666666; Value = 624790700000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Match the date-format/style already used for column A on data rows
    # (reuses the existing cell style rather than minting a new one).
    $ws.Range("A220").Copy()
    $ws.Range("A${row}").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.DateSerial
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:ROM2"
    $ws.Range("C${row}:F${row}").Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}
